# Adds all IG authors as contact
# The "Metadata" sheet already lists two "Contact" / "No display for
# ContactDetail" rows (rows 10-11). Two more IG authors need a "Contact"
# row each, so insert two rows right after the existing contact rows
# (i.e. before the "Jurisdiction" row) and fill them with the same
# Property/Value pair, duplicating the existing contact rows' formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Metadata")

# Insert two new rows (just the used A:B columns) before row 12
# ("Jurisdiction"), pushing the following rows down. Doing this twice,
# always at row 12, yields two fresh blank rows at 12 and 13.
$ws.Range("A12:B12").Insert()
$ws.Range("A12:B12").Insert()

# Match the formatting of the existing "Contact" row (row 11) by copying
# it onto the two freshly inserted rows.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A11:B11").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)

# Fill in the duplicated contact values (two more IG authors as contacts).
$ws.Cells.Item(12, 1).Value = "Contact"
$ws.Cells.Item(12, 2).Value = "No display for ContactDetail"
$ws.Cells.Item(13, 1).Value = "Contact"
$ws.Cells.Item(13, 2).Value = "No display for ContactDetail"
